# Design Document edit: correct "Al Al-Mohaiminul Islam Khan" to
# "Al-Mohaiminul Islam Khan" in the team roster table, and leave Word's
# "_GoBack" last-edit bookmark at that location (moving it away from its
# previous spot next to "Crop Attribute", where it collapses back into a
# single run).

$d = $word.ActiveDocument

# 1) Fix the duplicated "Al " prefix: "Al Al-" -> "Al-"
$d.Content.Find.Execute("Al Al-", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Al-", 2)

# 2) Re-locate the "_GoBack" bookmark to the start of the text we just
#    edited (this is what Word itself does after an edit). Adding a
#    bookmark with a name that already exists elsewhere in the document
#    implicitly removes the old one, since bookmark names must be unique.
$r = $d.Content
$r.Find.Execute("Al-", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$target = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $target)

# 3) Where "_GoBack" used to sit (after "Crop Attribute"), the paragraph
#    is now left with two adjacent runs ("Crop Attribute" + " "); collapse
#    them back into the single run Word produces once the bookmark is gone.
$d.Content.Find.Execute("Crop Attribute ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Crop Attribute ", 2)
